# Updated cryptos list with refreshed price and volume(1h) data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.786.37"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "2.223.63"
$ws.Range("E3").Value = "  -4.85%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'298.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.37%  "

$ws.Range("D6").Value = "'84.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.03%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.78%  "

$ws.Range("D10").Value = "'0.0781"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.04%  "

$ws.Range("D11").Value = "'29.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("D12").Value = "'46.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.96%  "

$ws.Range("E13").Value = "  -2.34%  "

$ws.Range("D14").Value = "2.575.73"
$ws.Range("E14").Value = "  -4.49%  "

$ws.Range("D15").Value = "'6.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.14%  "

$ws.Range("D16").Value = "'14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.38%  "

$ws.Range("D17").Value = "2.224.86"
$ws.Range("E17").Value = "  -5.38%  "

$ws.Range("D18").Value = "'0.717"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.39%  "

$ws.Range("D19").Value = "39.712.55"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("D20").Value = "0.0₃0876"
$ws.Range("E20").Value = "  -2.83%  "

$ws.Range("D21").Value = "'5.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.80%  "

$ws.Range("D22").Value = "'65.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.32%  "

$ws.Range("D23").Value = "'10.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.42%  "

$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("E26").Value = "  -5.04%  "

$ws.Range("D27").Value = "'1.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("D28").Value = "'22.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").Value = "'9.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").Value = "'32.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.45%  "

$ws.Range("D32").Value = "'149.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.69%  "

$ws.Range("D34").Value = "'4.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.48%  "

$ws.Range("E35").Value = "  -1.58%  "

$ws.Range("D36").Value = "'0.0701"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.33%  "

$ws.Range("D37").Value = "'16.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.20%  "

$ws.Range("E38").Value = "  -2.79%  "

$ws.Range("D39").Value = "'0.0975"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").Value = "'2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("D41").Value = "'1.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.91%  "

$ws.Range("D42").Value = "'3.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.25%  "

$ws.Range("D43").Value = "1.930.25"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("D45").Value = "'0.0265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("D47").Value = "'16.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.82%  "

$ws.Range("D48").Value = "'2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "

$ws.Range("D49").Value = "2.446.08"
$ws.Range("E49").Value = "  -4.24%  "

$ws.Range("D50").Value = "'71.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "

$ws.Range("D51").Value = "'88.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.60%  "
